# "Update to DB Ready"
# 1) Convert the Quantity column (E2:E30) from text-numbers to real numbers.
# 2) Append two new product rows (31 and 32) with the matching data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Column E (Quantity) rows 2-30: store as genuine numeric values instead
#    of text so downstream DB/number processing works correctly.
# ---------------------------------------------------------------------------
$quantities = @{
    2  = 4
    3  = 1
    4  = 2
    5  = 2
    6  = 2
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 3
    16 = 2
    17 = 1
    18 = 5
    19 = 1
    20 = 1
    21 = 3
    22 = 2
    23 = 2
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 2
    30 = 2
}

foreach ($r in $quantities.Keys) {
    $ws.Cells.Item($r, 5).Value = $quantities[$r]
}

# ---------------------------------------------------------------------------
# 2) Append the two new rows of product data.
# ---------------------------------------------------------------------------
# Columns whose values look like plain numbers but must stay text (matching
# the rest of the sheet, e.g. leading-zero item numbers) get a leading
# apostrophe so Excel keeps them as text instead of coercing to a number.
$numericLookingCols = @(1, 5, 9, 10, 11, 12)

$newRows = @(
    @{
        Row  = 31
        Data = @{
            1  = "000030"
            2  = "PILOT IGNITION"
            3  = ""
            4  = ""
            5  = "2"
            6  = "--"
            7  = "IGNITION"
            8  = "--"
            9  = "1"
            10 = "2"
            11 = "5"
            12 = "2"
            13 = "H-302"
            14 = "PILOT IGNITION"
            15 = "BTX"
            16 = "IG-ROD-PILOT"
            17 = "--"
        }
    },
    @{
        Row  = 32
        Data = @{
            1  = "000031"
            2  = "IGNITION TRANSFORMER"
            3  = "SIEMENS"
            4  = "ZA 20 100 LH 21"
            5  = "2"
            6  = "--"
            7  = "TRANSFORMER"
            8  = "220VAC"
            9  = "1"
            10 = "2"
            11 = "4"
            12 = "2"
            13 = "BOILER"
            14 = ""
            15 = "BOILER"
            16 = "IG-XF-SIEMENS-ZA 20 100 LH 21"
            17 = ""
        }
    }
)

foreach ($rowInfo in $newRows) {
    $r = $rowInfo.Row
    $data = $rowInfo.Data
    foreach ($c in $data.Keys) {
        $val = $data[$c]
        if ($val -ne "") {
            if ($numericLookingCols -contains $c) {
                $ws.Cells.Item($r, $c).Value = "'" + $val
            } else {
                $ws.Cells.Item($r, $c).Value = $val
            }
        }
    }
}
